$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 16:17"

# Row 4
$ws.Range("B4").Value = 6051431
$ws.Range("C4").Value = 4797
$ws.Range("D4").Value = 3348934
$ws.Range("E4").Value = 2517530
$ws.Range("G4").Value = 171
$ws.Range("H4").Value = 184967

# Row 6
$ws.Range("B6").Value = 3403555
$ws.Range("C6").Value = 18980
$ws.Range("D6").Value = 2596273
$ws.Range("E6").Value = 745425
$ws.Range("G6").Value = 163
$ws.Range("H6").Value = 61857

# Row 14
$ws.Range("D14").Value = 280165
$ws.Range("E14").Value = 91998
$ws.Range("G14").Value = 79
$ws.Range("H14").Value = 8129

# Row 23
$ws.Range("B23").Value = 241257
$ws.Range("C23").Value = 692
$ws.Range("E23").Value = 16402
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 9360

# Row 41
$ws.Range("A41").Value = "Kuwait"
$ws.Range("B41").Value = 83578
$ws.Range("C41").Value = 633
$ws.Range("D41").Value = 75320
$ws.Range("E41").Value = 7733
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 525

# Row 42
$ws.Range("A42").Value = "Belgica"
$ws.Range("B42").Value = 83500
$ws.Range("C42").Value = 470
$ws.Range("D42").Value = 18360
$ws.Range("E42").Value = 55256
$ws.Range("G42").Value = 5
$ws.Range("H42").Value = 9884

# Row 51
$ws.Range("B51").Value = 57074
$ws.Range("C51").Value = 401
$ws.Range("D51").Value = 41556
$ws.Range("E51").Value = 13703
$ws.Range("G51").Value = 6
$ws.Range("H51").Value = 1815

# Row 68
$ws.Range("B68").Value = 33630
$ws.Range("C68").Value = 241
$ws.Range("D68").Value = 19434
$ws.Range("E68").Value = 13629

# Row 69
$ws.Range("B69").Value = 31207
$ws.Range("C69").Value = 108
$ws.Range("D69").Value = 29802
$ws.Range("E69").Value = 696
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 709

# Row 85
$ws.Range("B85").Value = 14163
$ws.Range("C85").Value = 159
$ws.Range("D85").Value = 10933
$ws.Range("E85").Value = 2640
$ws.Range("G85").Value = 5
$ws.Range("H85").Value = 590

# Row 89
$ws.Range("B89").Value = 11779
$ws.Range("C89").Value = 178
$ws.Range("D89").Value = 10945
$ws.Range("E89").Value = 551
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 283

# Row 98
$ws.Range("A98").Value = "Tayikistan"
$ws.Range("B98").Value = 8481
$ws.Range("C98").Value = 32
$ws.Range("D98").Value = 7276
$ws.Range("E98").Value = 1137
$ws.Range("H98").Value = 68

# Row 99
$ws.Range("A99").Value = "Gabon"
$ws.Range("B99").Value = 8468
$ws.Range("D99").Value = 7066
$ws.Range("E99").Value = 1349
$ws.Range("H99").Value = 53

# Row 130
$ws.Range("A130").Value = "Uganda"
$ws.Range("B130").Value = 2756
$ws.Range("C130").Value = 77
$ws.Range("D130").Value = 1288
$ws.Range("E130").Value = 1440
$ws.Range("H130").Value = 28

# Row 131
$ws.Range("A131").Value = "Gambia"
$ws.Range("B131").Value = 2743
$ws.Range("D131").Value = 638
$ws.Range("E131").Value = 2012
$ws.Range("H131").Value = 93

# Row 132
$ws.Range("A132").Value = "Mali"
$ws.Range("B132").Value = 2730
$ws.Range("D132").Value = 2054
$ws.Range("E132").Value = 550
$ws.Range("H132").Value = 126

# Row 139
$ws.Range("B139").Value = 2092
$ws.Range("C139").Value = 5
$ws.Range("E139").Value = 118

# Row 184
$ws.Range("A184").Value = "Gibraltar"
$ws.Range("B184").Value = 274
$ws.Range("C184").Value = 2
$ws.Range("D184").Value = 221
$ws.Range("E184").Value = 53

# Row 185
$ws.Range("A185").Value = "Camboya"
$ws.Range("B185").Value = 273
$ws.Range("D185").Value = 265
$ws.Range("E185").Value = 8
